# Update header labels on every worksheet:
#   "% Not falling within the PiN dimensions" -> "% Not in need"
#   "# Not falling within the PiN dimensions" -> "# Not in need"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    foreach ($cell in $ws.UsedRange.Cells) {
        $v = $cell.Value2
        if ($v -eq "% Not falling within the PiN dimensions") {
            $cell.Value = "% Not in need"
        }
        elseif ($v -eq "# Not falling within the PiN dimensions") {
            $cell.Value = "# Not in need"
        }
    }
}
